$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.238.06"
$ws.Range("E2").Value = "  +4.62%  "
$ws.Range("D3").Value = "3.241.67"
$ws.Range("E3").Value = "  +2.59%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'576.00"
$ws.Range("E5").Value = "  +2.14%  "
$ws.Range("D6").Value = "'180.12"
$ws.Range("E6").Value = "  +5.94%  "
$ws.Range("D8").Value = "'0.591"
$ws.Range("E8").Value = "  -2.96%  "
$ws.Range("D9").Value = "3.237.03"
$ws.Range("E9").Value = "  +2.56%  "
$ws.Range("E10").Value = "  +4.30%  "
$ws.Range("D11").Value = "'6.78"
$ws.Range("E11").Value = "  +3.38%  "
$ws.Range("D12").Value = "'0.411"
$ws.Range("E12").Value = "  +4.57%  "
$ws.Range("D13").Value = "3.800.12"
$ws.Range("E13").Value = "  +2.58%  "
$ws.Range("E14").Value = "  +1.21%  "
$ws.Range("D15").Value = "'27.86"
$ws.Range("E15").Value = "  +2.35%  "
$ws.Range("D16").Value = "67.217.50"
$ws.Range("E16").Value = "  +4.78%  "
$ws.Range("D17").Value = "'0.0000167"
$ws.Range("E17").Value = "  +2.74%  "
$ws.Range("D18").Value = "3.253.50"
$ws.Range("E18").Value = "  +3.44%  "
$ws.Range("D19").Value = "'5.80"
$ws.Range("E19").Value = "  +1.27%  "
$ws.Range("D20").Value = "'13.37"
$ws.Range("E20").Value = "  +3.48%  "
$ws.Range("D21").Value = "'372.65"
$ws.Range("E21").Value = "  +5.34%  "
$ws.Range("D22").Value = "'7.57"
$ws.Range("E22").Value = "  +5.15%  "
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").Value = "'70.84"
$ws.Range("E24").Value = "  +4.12%  "
$ws.Range("D25").Value = "'0.509"
$ws.Range("E25").Value = "  +1.76%  "
$ws.Range("E26").Value = "  +1.96%  "
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("D28").Value = "'0.180"
$ws.Range("E28").Value = "  +3.26%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.38%  "
$ws.Range("E30").Value = "  +4.66%  "
$ws.Range("D31").Value = "'5.66"
$ws.Range("E31").Value = "  +4.06%  "
$ws.Range("D32").Value = "'22.57"
$ws.Range("E32").Value = "  +2.95%  "
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("E34").Value = "  +5.05%  "
$ws.Range("D35").Value = "'6.85"
$ws.Range("E35").Value = "  +3.39%  "
$ws.Range("D36").Value = "'1.49"
$ws.Range("E36").Value = "  +4.44%  "
$ws.Range("D37").Value = "'161.82"
$ws.Range("E37").Value = "  +5.30%  "
$ws.Range("D38").Value = "'0.853"
$ws.Range("E38").Value = "  +3.98%  "
$ws.Range("D39").Value = "'1.84"
$ws.Range("E39").Value = "  +8.98%  "
$ws.Range("D40").Value = "'6.80"
$ws.Range("E40").Value = "  +13.79%  "
$ws.Range("D41").Value = "'26.67"
$ws.Range("E41").Value = "  +2.03%  "
$ws.Range("E42").Value = "  +6.07%  "
$ws.Range("D43").Value = "'360.94"
$ws.Range("E43").Value = "  +14.18%  "
$ws.Range("D44").Value = "'4.38"
$ws.Range("E44").Value = "  +5.22%  "
$ws.Range("D45").Value = "2.697.28"
$ws.Range("E45").Value = "  +2.11%  "
$ws.Range("D46").Value = "'25.60"
$ws.Range("E46").Value = "  +6.80%  "
$ws.Range("D47").Value = "'40.50"
$ws.Range("E47").Value = "  +2.92%  "
$ws.Range("D48").Value = "'0.0672"
$ws.Range("E48").Value = "  +3.26%  "
$ws.Range("D49").Value = "'0.0277"
$ws.Range("E49").Value = "  +2.49%  "
$ws.Range("D50").Value = "'0.995"
$ws.Range("E50").Value = "  +6.40%  "
$ws.Range("D51").Value = "'0.102"
$ws.Range("E51").Value = "  +1.60%  "
